$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = -8.667399999999994
    14  = -8.2547
    21  = -8.2699
    23  = -8.047899999999997
    25  = -7.839599999999994
    26  = -7.226700000000005
    29  = -7.378400000000003
    53  = -8.268299999999996
    57  = -8.140699999999995
    59  = -8.291099999999993
    69  = -7.689899999999994
    79  = -7.381500000000004
    83  = -9.091499999999998
    91  = -8.035499999999997
    93  = -6.908499999999991
    103 = -7.802799999999995
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
